# Update column F (dSF) values for rows 2-26 per the repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = 1
    4  = 4
    5  = -2
    6  = 3
    7  = 2
    8  = 3
    10 = -4
    11 = 7
    12 = 7
    13 = 7
    14 = 2
    15 = 1
    16 = -2
    17 = 2
    19 = -2
    20 = 8
    22 = -4
    23 = 2
    24 = 1
    25 = -1
    26 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
